# ajuste: corrigindo as categorias
#
# Adds two new header columns (S = "Idade ignorada", T = "Total"),
# adds two new data rows ("Outros" and "Total"), and fills in the
# S/T values for the pre-existing category rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New header cells (row 1) ---
$ws.Range("S1").Value = "Idade ignorada"
$ws.Range("T1").Value = "Total"

# --- New column T (Total) for the existing category rows (2-6) ---
$ws.Range("T2").Value = 2086
$ws.Range("T3").Value = 256
$ws.Range("T4").Value = 709
$ws.Range("T5").Value = 378
$ws.Range("T6").Value = 1487

# --- New row 7: "Outros" ---
$ws.Range("A7").Value = "Outros"
$ws.Range("B7").Value = 114
$ws.Range("C7").Value = 3
$ws.Range("D7").Value = 6
$ws.Range("E7").Value = 45
$ws.Range("F7").Value = 75
$ws.Range("G7").Value = 77
$ws.Range("H7").Value = 78
$ws.Range("I7").Value = 105
$ws.Range("J7").Value = 127
$ws.Range("K7").Value = 148
$ws.Range("L7").Value = 204
$ws.Range("M7").Value = 224
$ws.Range("N7").Value = 287
$ws.Range("O7").Value = 342
$ws.Range("P7").Value = 346
$ws.Range("Q7").Value = 365
$ws.Range("R7").Value = 1126
$ws.Range("S7").Value = 1
$ws.Range("T7").Value = 3673

# --- New row 8: "Total" ---
$ws.Range("A8").Value = "Total"
$ws.Range("B8").Value = 123
$ws.Range("C8").Value = 5
$ws.Range("D8").Value = 6
$ws.Range("E8").Value = 52
$ws.Range("F8").Value = 88
$ws.Range("G8").Value = 100
$ws.Range("H8").Value = 119
$ws.Range("I8").Value = 166
$ws.Range("J8").Value = 212
$ws.Range("K8").Value = 294
$ws.Range("L8").Value = 418
$ws.Range("M8").Value = 540
$ws.Range("N8").Value = 748
$ws.Range("O8").Value = 893
$ws.Range("P8").Value = 892
$ws.Range("Q8").Value = 947
$ws.Range("R8").Value = 2985
$ws.Range("S8").Value = 1
$ws.Range("T8").Value = 8589
